$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C23").Value = "Final project presentations"
$ws.Range("E23").ClearContents()
$ws.Range("C24").Value = "Final project presentations; workshop final papers"
$ws.Range("E24").Value = "Component 4 draft"
$ws.Range("G19").Value = "Interpreting results"

$ws.Range("D21").Select()
